# Statistics workbook update — add newly-tracked click/nickname events for
# the "Анастасия Зуева" user (Telegram id 1318927560, nickname
# ZuevaAnastasiya) on the "Первый слой меню" and "Успешные клики" sheets,
# and widen the name / id columns so the longer values stay readable.

$wb = $excel.ActiveWorkbook

# ---- "Первый слой меню" sheet: 3 new button-click rows ----------------
$wsMenu = $wb.Worksheets.Item("Первый слой меню")

$wsMenu.Range("A20").Value = "Анастасия"
$wsMenu.Range("B20").Value = "Зуева"
$wsMenu.Range("C20").Value = "about"
$wsMenu.Range("D20").Value = "21.08.2024 16:38:46"

$wsMenu.Range("A21").Value = "Анастасия"
$wsMenu.Range("B21").Value = "Зуева"
$wsMenu.Range("C21").Value = "criteries"
$wsMenu.Range("D21").Value = "21.08.2024 16:39:08"

$wsMenu.Range("A22").Value = "Анастасия"
$wsMenu.Range("B22").Value = "Зуева"
$wsMenu.Range("C22").Value = "calendar"
$wsMenu.Range("D22").Value = "21.08.2024 16:39:47"

# Column A ("Имя") needs to be wider now that longer names are present.
$wsMenu.Columns.Item(1).ColumnWidth = 10.83

# ---- "Успешные клики" sheet: 2 new successful-click rows ---------------
$wsSuccess = $wb.Worksheets.Item("Успешные клики")

$wsSuccess.Range("A5").Value = "Анастасия"
$wsSuccess.Range("B5").Value = "Зуева"
$wsSuccess.Range("C5").Value = "ZuevaAnastasiya"
$wsSuccess.Range("D5").Value = "21.08.2024 16:39:12"

$wsSuccess.Range("A6").Value = "Анастасия"
$wsSuccess.Range("B6").Value = "Зуева"
$wsSuccess.Range("C6").Value = "ZuevaAnastasiya"
$wsSuccess.Range("D6").Value = "21.08.2024 16:39:12"

# Column A ("Имя") and column C ("ID в Телеграм" / nickname) both widen.
$wsSuccess.Columns.Item(1).ColumnWidth = 10.83
$wsSuccess.Columns.Item(3).ColumnWidth = 16.83
